$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 223.64
$ws.Cells.Item(15, 9).Value = 223.64
$ws.Cells.Item(15, 11).Value = 670.92
$ws.Cells.Item(15, 13).Value = -501.92
# Row 68
$ws.Cells.Item(68, 8).Value = 32000
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).ClearContents()
# Row 70
$ws.Cells.Item(70, 8).Value = 11977434
$ws.Cells.Item(70, 9).Value = 55889656
$ws.Cells.Item(70, 10).Value = 1372.7273
$ws.Cells.Item(70, 11).Value = 167668968
$ws.Cells.Item(70, 12).Value = 4118.1819
$ws.Cells.Item(70, 13).Value = -167668698
$ws.Cells.Item(70, 14).Value = -4658.1819
# Row 71
$ws.Cells.Item(71, 8).Value = 32000
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).ClearContents()
# Row 73
$ws.Cells.Item(73, 8).Value = 11977434
$ws.Cells.Item(73, 9).Value = 55889656
$ws.Cells.Item(73, 10).Value = 1372.7273
$ws.Cells.Item(73, 11).Value = 167668968
$ws.Cells.Item(73, 12).Value = 4118.1819
$ws.Cells.Item(73, 13).Value = -167668032
$ws.Cells.Item(73, 14).Value = -5990.1819
# Row 74
$ws.Cells.Item(74, 8).Value = 3780
$ws.Cells.Item(74, 9).Value = 3250
$ws.Cells.Item(74, 11).Value = 3250
$ws.Cells.Item(74, 13).Value = -2314
# Row 77
$ws.Cells.Item(77, 8).Value = 3780
$ws.Cells.Item(77, 9).Value = 3250
$ws.Cells.Item(77, 11).Value = 16250
$ws.Cells.Item(77, 13).Value = -11570
# Row 113
$ws.Cells.Item(113, 8).Value = 3600
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 3750
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 3750
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(113, 14).Value = -10258
# Row 129
$ws.Cells.Item(129, 8).Value = 1133.6061
$ws.Cells.Item(129, 10).Value = 1221.0741
$ws.Cells.Item(129, 12).Value = 3663.2223
$ws.Cells.Item(129, 14).Value = -13663.2223
# Row 134
$ws.Cells.Item(134, 8).Value = 111672.64
$ws.Cells.Item(134, 10).Value = 111672.64
$ws.Cells.Item(134, 12).Value = 111672.64
$ws.Cells.Item(134, 14).Value = -121812.64
# Row 138
$ws.Cells.Item(138, 8).Value = 2155498.8
$ws.Cells.Item(138, 9).Value = 4447556
$ws.Cells.Item(138, 10).Value = 6695.25
$ws.Cells.Item(138, 11).Value = 13342668
$ws.Cells.Item(138, 12).Value = 20085.75
$ws.Cells.Item(138, 13).Value = -13337528
$ws.Cells.Item(138, 14).Value = -30365.75

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Cells.Item(31, 8).Value = 6488.3
$ws.Cells.Item(31, 9).Value = 6488.3
$ws.Cells.Item(31, 11).Value = 6488.3
$ws.Cells.Item(31, 13).Value = -6194.3
# Row 32
$ws.Cells.Item(32, 8).Value = 64401.734
$ws.Cells.Item(32, 9).Value = 74501.5
$ws.Cells.Item(32, 10).Value = 52859.145
$ws.Cells.Item(32, 11).Value = 74501.5
$ws.Cells.Item(32, 12).Value = 52859.145
$ws.Cells.Item(32, 13).Value = -74214.5
$ws.Cells.Item(32, 14).Value = -53433.145
# Row 61
$ws.Cells.Item(61, 8).Value = 3286.2144
$ws.Cells.Item(61, 9).Value = 2886.5557
$ws.Cells.Item(61, 10).Value = 4005.6
$ws.Cells.Item(61, 11).Value = 2886.5557
$ws.Cells.Item(61, 12).Value = 4005.6
$ws.Cells.Item(61, 13).Value = -2674.5557
$ws.Cells.Item(61, 14).Value = -4429.6
# Row 74
$ws.Cells.Item(74, 8).Value = 1485.9
$ws.Cells.Item(74, 9).Value = 1479.12
$ws.Cells.Item(74, 11).Value = 1479.12
$ws.Cells.Item(74, 13).Value = -605.1199999999999
# Row 77
$ws.Cells.Item(77, 8).Value = 1485.9
$ws.Cells.Item(77, 9).Value = 1479.12
$ws.Cells.Item(77, 11).Value = 7395.599999999999
$ws.Cells.Item(77, 13).Value = -3027.599999999999
# Row 109
$ws.Cells.Item(109, 8).Value = 33000
$ws.Cells.Item(109, 10).Value = 33000
$ws.Cells.Item(109, 12).Value = 33000
$ws.Cells.Item(109, 14).Value = -35774
# Row 112
$ws.Cells.Item(112, 8).Value = 23598.572
$ws.Cells.Item(112, 10).Value = 23598.572
$ws.Cells.Item(112, 12).Value = 23598.572
$ws.Cells.Item(112, 14).Value = -26552.572
# Row 132
$ws.Cells.Item(132, 8).Value = 2711.9333
$ws.Cells.Item(132, 9).Value = 2390
$ws.Cells.Item(132, 11).Value = 7170
$ws.Cells.Item(132, 13).Value = -4640
# Row 134
$ws.Cells.Item(134, 8).Value = 52592.332
$ws.Cells.Item(134, 10).Value = 52592.332
$ws.Cells.Item(134, 12).Value = 52592.332
$ws.Cells.Item(134, 14).Value = -62732.332
# Row 136
$ws.Cells.Item(136, 8).Value = 3286.2144
$ws.Cells.Item(136, 9).Value = 2886.5557
$ws.Cells.Item(136, 10).Value = 4005.6
$ws.Cells.Item(136, 11).Value = 8659.667099999999
$ws.Cells.Item(136, 12).Value = 12016.8
$ws.Cells.Item(136, 13).Value = -6109.667099999999
$ws.Cells.Item(136, 14).Value = -17116.8
# Row 138
$ws.Cells.Item(138, 8).Value = 62850
$ws.Cells.Item(138, 10).Value = 62850
$ws.Cells.Item(138, 12).Value = 62850
$ws.Cells.Item(138, 14).Value = -73130

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Cells.Item(81, 8).Value = 16906
$ws.Cells.Item(81, 10).Value = 16906
$ws.Cells.Item(81, 12).Value = 16906
$ws.Cells.Item(81, 14).Value = -19028
# Row 84
$ws.Cells.Item(84, 8).Value = 16906
$ws.Cells.Item(84, 10).Value = 16906
$ws.Cells.Item(84, 12).Value = 50718
$ws.Cells.Item(84, 14).Value = -61326
# Row 86
$ws.Cells.Item(86, 8).Value = 236353.44
$ws.Cells.Item(86, 9).Value = 23039.4
$ws.Cells.Item(86, 10).Value = 502996
$ws.Cells.Item(86, 11).Value = 23039.4
$ws.Cells.Item(86, 12).Value = 502996
$ws.Cells.Item(86, 13).Value = -21916.4
$ws.Cells.Item(86, 14).Value = -505242
# Row 89
$ws.Cells.Item(89, 8).Value = 236353.44
$ws.Cells.Item(89, 9).Value = 23039.4
$ws.Cells.Item(89, 10).Value = 502996
$ws.Cells.Item(89, 11).Value = 115197
$ws.Cells.Item(89, 12).Value = 2514980
$ws.Cells.Item(89, 13).Value = -109581
$ws.Cells.Item(89, 14).Value = -2526212
# Row 107
$ws.Cells.Item(107, 8).Value = 20942.037
$ws.Cells.Item(107, 9).Value = 27321.15
$ws.Cells.Item(107, 10).Value = 2716
$ws.Cells.Item(107, 11).Value = 27321.15
$ws.Cells.Item(107, 12).Value = 2716
$ws.Cells.Item(107, 13).Value = -25401.15
$ws.Cells.Item(107, 14).Value = -6556
# Row 134
$ws.Cells.Item(134, 8).Value = 2721.6416
$ws.Cells.Item(134, 9).Value = 1918.4857
$ws.Cells.Item(134, 10).Value = 4283.3335
$ws.Cells.Item(134, 11).Value = 5755.4571
$ws.Cells.Item(134, 12).Value = 12850.0005
$ws.Cells.Item(134, 13).Value = -3220.4571
$ws.Cells.Item(134, 14).Value = -17920.0005

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 3719.5518
$ws.Cells.Item(22, 9).Value = 4564.1304
$ws.Cells.Item(22, 10).Value = 482
$ws.Cells.Item(22, 11).Value = 4564.1304
$ws.Cells.Item(22, 12).Value = 482
$ws.Cells.Item(22, 13).Value = -4214.1304
$ws.Cells.Item(22, 14).Value = -1182
# Row 99
$ws.Cells.Item(99, 8).Value = 4154.5454
$ws.Cells.Item(99, 9).Value = 4271.4287
$ws.Cells.Item(99, 11).Value = 4271.4287
$ws.Cells.Item(99, 13).Value = -2773.4287
# Row 126
$ws.Cells.Item(126, 8).Value = 4154.5454
$ws.Cells.Item(126, 9).Value = 4271.4287
$ws.Cells.Item(126, 11).Value = 12814.2861
$ws.Cells.Item(126, 13).Value = -10344.2861
# Row 133
$ws.Cells.Item(133, 8).Value = 49610
$ws.Cells.Item(133, 10).Value = 49610
$ws.Cells.Item(133, 12).Value = 49610
$ws.Cells.Item(133, 14).Value = -54670
# Row 134
$ws.Cells.Item(134, 8).Value = 2035.375
$ws.Cells.Item(134, 9).Value = 1996.2858
$ws.Cells.Item(134, 10).Value = 2090.1
$ws.Cells.Item(134, 11).Value = 5988.857400000001
$ws.Cells.Item(134, 12).Value = 6270.299999999999
$ws.Cells.Item(134, 13).Value = -3453.857400000001
$ws.Cells.Item(134, 14).Value = -11340.3
# Row 137
$ws.Cells.Item(137, 8).Value = 74465
$ws.Cells.Item(137, 10).Value = 74465
$ws.Cells.Item(137, 12).Value = 74465
$ws.Cells.Item(137, 14).Value = -84665
# Row 140
$ws.Cells.Item(140, 8).Value = 78210
$ws.Cells.Item(140, 10).Value = 78210
$ws.Cells.Item(140, 12).Value = 78210
$ws.Cells.Item(140, 14).Value = -88570

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 133
$ws.Cells.Item(133, 8).Value = 2905.2222
$ws.Cells.Item(133, 9).Value = 1606.75
$ws.Cells.Item(133, 10).Value = 3944
$ws.Cells.Item(133, 11).Value = 4820.25
$ws.Cells.Item(133, 12).Value = 11832
$ws.Cells.Item(133, 13).Value = 239.75
$ws.Cells.Item(133, 14).Value = -21952

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Cells.Item(41, 8).Value = 5262.75
$ws.Cells.Item(41, 9).Value = 1025.5
$ws.Cells.Item(41, 11).Value = 1025.5
$ws.Cells.Item(41, 13).Value = -670.5
# Row 111
$ws.Cells.Item(111, 8).Value = 30000
$ws.Cells.Item(111, 10).Value = 30000
$ws.Cells.Item(111, 12).Value = 30000
$ws.Cells.Item(111, 14).Value = -36134
# Row 122
$ws.Cells.Item(122, 8).Value = 5387.7646
$ws.Cells.Item(122, 9).Value = 5132.6665
$ws.Cells.Item(122, 10).Value = 6000
$ws.Cells.Item(122, 11).Value = 15397.9995
$ws.Cells.Item(122, 12).Value = 18000
$ws.Cells.Item(122, 13).Value = -12947.9995
$ws.Cells.Item(122, 14).Value = -22900
# Row 132
$ws.Cells.Item(132, 8).Value = 1850.4
$ws.Cells.Item(132, 9).Value = 1500.5555
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 4501.666499999999
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -1971.666499999999
$ws.Cells.Item(132, 14).Value = -20057

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 110
$ws.Cells.Item(110, 8).Value = 39000
$ws.Cells.Item(110, 10).Value = 39000
$ws.Cells.Item(110, 12).Value = 39000
$ws.Cells.Item(110, 14).Value = -47180
# Row 132
$ws.Cells.Item(132, 8).Value = 4570.5
$ws.Cells.Item(132, 9).Value = 4243.7144
$ws.Cells.Item(132, 10).Value = 5333
$ws.Cells.Item(132, 11).Value = 12731.1432
$ws.Cells.Item(132, 12).Value = 15999
$ws.Cells.Item(132, 13).Value = -10201.1432
$ws.Cells.Item(132, 14).Value = -21059
# Row 133
$ws.Cells.Item(133, 8).Value = 60956.76
$ws.Cells.Item(133, 10).Value = 60956.76
$ws.Cells.Item(133, 12).Value = 60956.76
$ws.Cells.Item(133, 14).Value = -66016.76000000001
# Row 136
$ws.Cells.Item(136, 8).Value = 8726.25
$ws.Cells.Item(136, 9).Value = 10167.667
$ws.Cells.Item(136, 10).Value = 4402
$ws.Cells.Item(136, 11).Value = 30503.001
$ws.Cells.Item(136, 12).Value = 13206
$ws.Cells.Item(136, 13).Value = -27953.001
$ws.Cells.Item(136, 14).Value = -18306
# Row 140
$ws.Cells.Item(140, 8).Value = 83350
$ws.Cells.Item(140, 10).Value = 83350
$ws.Cells.Item(140, 12).Value = 83350
$ws.Cells.Item(140, 14).Value = -93710
# Row 141
$ws.Cells.Item(141, 8).Value = 47000
$ws.Cells.Item(141, 10).Value = 47000
$ws.Cells.Item(141, 12).Value = 47000
$ws.Cells.Item(141, 14).Value = -57360

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 4000
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).Value = -5248
# Row 65
$ws.Cells.Item(65, 8).Value = 4000
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).Value = -26240
# Row 122
$ws.Cells.Item(122, 8).Value = 50002260
$ws.Cells.Item(122, 9).Value = 125001000
$ws.Cells.Item(122, 11).Value = 375003000
$ws.Cells.Item(122, 13).Value = -375000550
# Row 138
$ws.Cells.Item(138, 8).Value = 109099
$ws.Cells.Item(138, 10).Value = 109099
$ws.Cells.Item(138, 12).Value = 109099
$ws.Cells.Item(138, 14).Value = -119379
# Row 140
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
